# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to match the refreshed output generated at 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 1781
    4  = 1677
    5  = 398
    10 = 244
    11 = 24
    12 = 90
    14 = 242
    16 = 35
    18 = 71
    19 = 242
    20 = 35
    21 = 441
    22 = 341
    27 = 747
    28 = 2553
    31 = 509
    32 = 829
    34 = 444
    35 = 248
    37 = 442
    38 = 576
    39 = 420
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
